$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 39.27015933333333
$ws.Range("H2").Value = 117.810478
$ws.Range("I2").Value = 0.2257020976862494
$ws.Range("J2").Value = 0.2257020976862494
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 117.044563
$ws.Range("N2").Value = 351.133689
$ws.Range("O2").Value = 0.3245365645427815
$ws.Range("P2").Value = 0.3245365645427815
$ws.Range("Q2").Value = 4596.358638110371
$ws.Range("R2").Value = 41367.22774299334
$ws.Range("S2").Value = 0.07324858339319465
$ws.Range("T2").Value = 0.07324858339319465
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 39.27015933333333
$ws.Range("H3").Value = 117.810478
$ws.Range("I3").Value = 0.2257020976862494
$ws.Range("J3").Value = 0.2257020976862494
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 101.5800373333333
$ws.Range("N3").Value = 304.740112
$ws.Range("O3").Value = 0.281657135515876
$ws.Range("P3").Value = 0.281657135515876
$ws.Range("Q3").Value = 3989.064251165948
$ws.Range("R3").Value = 35901.57826049354
$ws.Range("S3").Value = 0.06357060631423345
$ws.Range("T3").Value = 0.06357060631423343
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 39.27015933333333
$ws.Range("H4").Value = 117.810478
$ws.Range("I4").Value = 0.2257020976862494
$ws.Range("J4").Value = 0.2257020976862494
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 142.0267893333333
$ws.Range("N4").Value = 426.080368
$ws.Range("O4").Value = 0.3938062999413425
$ws.Range("P4").Value = 0.3938062999413425
$ws.Range("Q4").Value = 5577.414646721767
$ws.Range("R4").Value = 50196.7318204959
$ws.Range("S4").Value = 0.08888290797882134
$ws.Range("T4").Value = 0.08888290797882134
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 119.3024773333333
$ws.Range("H5").Value = 357.907432
$ws.Range("I5").Value = 0.6856814398113102
$ws.Range("J5").Value = 0.6856814398113102
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 117.044563
$ws.Range("N5").Value = 351.133689
$ws.Range("O5").Value = 0.3245365645427815
$ws.Range("P5").Value = 0.3245365645427815
$ws.Range("Q5").Value = 13963.70632429741
$ws.Range("R5").Value = 125673.3569186767
$ws.Range("S5").Value = 0.2225286988471106
$ws.Range("T5").Value = 0.2225286988471106
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 119.3024773333333
$ws.Range("H6").Value = 357.907432
$ws.Range("I6").Value = 0.6856814398113102
$ws.Range("J6").Value = 0.6856814398113102
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 101.5800373333333
$ws.Range("N6").Value = 304.740112
$ws.Range("O6").Value = 0.281657135515876
$ws.Range("P6").Value = 0.281657135515876
$ws.Range("Q6").Value = 12118.75010147916
$ws.Range("R6").Value = 109068.7509133124
$ws.Range("S6").Value = 0.1931270702136552
$ws.Range("T6").Value = 0.1931270702136552
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 119.3024773333333
$ws.Range("H7").Value = 357.907432
$ws.Range("I7").Value = 0.6856814398113102
$ws.Range("J7").Value = 0.6856814398113102
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 142.0267893333333
$ws.Range("N7").Value = 426.080368
$ws.Range("O7").Value = 0.3938062999413425
$ws.Range("P7").Value = 0.3938062999413425
$ws.Range("Q7").Value = 16944.14781516611
$ws.Range("R7").Value = 152497.330336495
$ws.Range("S7").Value = 0.2700256707505445
$ws.Range("T7").Value = 0.2700256707505445
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 15.418477
$ws.Range("H8").Value = 46.255431
$ws.Range("I8").Value = 0.08861646250244033
$ws.Range("J8").Value = 0.08861646250244033
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 117.044563
$ws.Range("N8").Value = 351.133689
$ws.Range("O8").Value = 0.3245365645427815
$ws.Range("P8").Value = 0.3245365645427815
$ws.Range("Q8").Value = 1804.648902590551
$ws.Range("R8").Value = 16241.84012331496
$ws.Range("S8").Value = 0.0287592823024762
$ws.Range("T8").Value = 0.0287592823024762
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 15.418477
$ws.Range("H9").Value = 46.255431
$ws.Range("I9").Value = 0.08861646250244033
$ws.Range("J9").Value = 0.08861646250244033
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 101.5800373333333
$ws.Range("N9").Value = 304.740112
$ws.Range("O9").Value = 0.281657135515876
$ws.Range("P9").Value = 0.281657135515876
$ws.Range("Q9").Value = 1566.209469283142
$ws.Range("R9").Value = 14095.88522354827
$ws.Range("S9").Value = 0.02495945898798738
$ws.Range("T9").Value = 0.02495945898798738
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 15.418477
$ws.Range("H10").Value = 46.255431
$ws.Range("I10").Value = 0.08861646250244033
$ws.Range("J10").Value = 0.08861646250244033
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 142.0267893333333
$ws.Range("N10").Value = 426.080368
$ws.Range("O10").Value = 0.3938062999413425
$ws.Range("P10").Value = 0.3938062999413425
$ws.Range("Q10").Value = 2189.836784719846
$ws.Range("R10").Value = 19708.53106247861
$ws.Range("S10").Value = 0.03489772121197675
$ws.Range("T10").Value = 0.03489772121197675
